# Update MicroUI diagrams in the VDG (replace usage of the 'platform' terminology)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the auto "Last updated" date field that PowerPoint stamps on the
#    slide master + every layout (re-saving the deck recomputes the cached
#    display text of the datetimeFigureOut field).
# ---------------------------------------------------------------------------
$newDate = "01/12/2025"

$master = $p.SlideMaster
$dateShape = $master.Shapes.Item("Date Placeholder 3")
$dateShape.TextFrame.TextRange.Text = $newDate

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 1 - MicroUI diagram relabelling / relayout.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# "Provided by user" -> "Provided by user(s)" (textbox grows a little wider)
$shpUser = $s.Shapes.Item("TextBox 76")
$shpUser.TextFrame.TextRange.Text = "Provided by user(s)"
$shpUser.Width = 111.99425196850393

# "Provided by platform" -> "Provided by VEE Port" (textbox shifts right)
$shpVeePort = $s.Shapes.Item("TextBox 77")
$shpVeePort.TextFrame.TextRange.Text = "Provided by VEE Port"
$shpVeePort.Left = 312.3668503937008

# "Platform" -> "Embedded" (textbox shifts left and grows wider)
$shpEmbedded = $s.Shapes.Item("TextBox 78")
$shpEmbedded.TextFrame.TextRange.Text = "Embedded"
$shpEmbedded.Left = 203.29433070866142
$shpEmbedded.Width = 74.70795275590551

# The small coloured marker dot next to "Provided by VEE Port" moves in step
# with its label. There are two shapes sharing the name "Rounded Rectangle
# 164" - the one we need is the second one (index 28) at x=3714317 EMU.
$shpDot = $s.Shapes.Item(28)
$shpDot.Left = 297.4072440944882
